$d = $word.ActiveDocument

# Locate the paragraph that reads exactly "File ExcessaoEmailInvalido.php"
# (there is a similarly named paragraph "File Utilidades_ExcessaoEmailInvalido.php.html"
# further down, so we match the full paragraph text incl. the trailing paragraph mark to
# make sure we grab the right one).
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "File ExcessaoEmailInvalido.php`r") {
        $p1 = $p
        break
    }
}

if ($p1 -eq $null) {
    throw "Could not find paragraph 'File ExcessaoEmailInvalido.php'"
}

# The next paragraph is the (empty) one holding the _GoBack bookmark, and the one
# after that is "File ExcessaoNomeInvalido.php".
$p2 = $p1.Next()
$p3 = $p2.Next()

if ($p3.Range.Text -ne "File ExcessaoNomeInvalido.php`r") {
    throw "Unexpected document structure near ExcessaoEmailInvalido paragraph"
}

# Replace the three paragraphs (File Email / empty+bookmark / File Nome) with:
#   File ExcessaoEmailInvalido.php<space>
#   Classe ExcessaoEmailInvalido           (new paragraph, numId 26)
#   <bookmark _GoBack>File ExcessaoNomeInvalido.php   (merged into the former empty paragraph)
$rng = $d.Range($p1.Range.Start, $p3.Range.End)

$fragment = ""
$fragment += '<w:p w:rsidR="008A40DF" w:rsidRDefault="008A40DF" w:rsidP="008A40DF"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr>'
$fragment += '<w:r><w:t xml:space="preserve">File </w:t></w:r>'
$fragment += '<w:proofErr w:type="gramStart"/><w:r w:rsidRPr="008A40DF"><w:t>ExcessaoEmailInvalido</w:t></w:r><w:proofErr w:type="gramEnd"/>'
$fragment += '<w:r><w:t>.php</w:t></w:r>'
$fragment += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$fragment += '</w:p>'
$fragment += '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr></w:pPr>'
$fragment += '<w:r><w:t xml:space="preserve">Classe </w:t></w:r>'
$fragment += '<w:proofErr w:type="gramStart"/><w:r><w:t>ExcessaoEmailInvalido</w:t></w:r><w:proofErr w:type="gramEnd"/>'
$fragment += '</w:p>'
$fragment += '<w:p w:rsidR="008B1E1B" w:rsidRDefault="008B1E1B" w:rsidP="008A40DF"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr>'
$fragment += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$fragment += '<w:r><w:t xml:space="preserve">File </w:t></w:r>'
$fragment += '<w:r><w:t>ExcessaoNomeInvalido</w:t></w:r>'
$fragment += '<w:r><w:t>.php</w:t></w:r>'
$fragment += '</w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$rng.InsertXML($xml)
